# Replace the placeholder "Business description..." text in the CI_Descrip /
# Dependency_Descrip columns (C and G) for the data rows (2-42) with the full
# business-description text (the same text already used in row 43), and grow
# those rows so the five-line description is fully visible. Also move the
# sheet's active selection to H3 (and drop the stale "scrolled to row 40"
# top-left anchor that came with the old D46 selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$longDescription = "Business description...`nAttribute 1: Value`nAttribute 2: Value`nAttribute 3: Value`n "

for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 3).Value = $longDescription
    $ws.Cells.Item($row, 7).Value = $longDescription
    $ws.Rows.Item($row).RowHeight = 52.5
}

$ws.Range("H3").Select()
